$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $value into $cellRef as literal TEXT, even when the string
# looks numeric (e.g. "53.07" or "0.967"). Assigning .Value directly would
# let Excel's COM layer auto-coerce such strings into floating point
# numbers, which corrupts values like "63.454.68" style prices and loses
# trailing zeros (e.g. "549.16" -> 549.16 displayed as "549.16" is fine, but
# "53.10" -> 53.1). Forcing NumberFormat "@" (Text) before the assignment
# keeps the literal text, and resetting the Style back to "Normal"
# afterwards avoids leaving a stray number-format style on the cell.
function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '63.454.68'
Set-TextValue 'E2' '  -3.13%  '
Set-TextValue 'D3' '3.315.09'
Set-TextValue 'E3' '  -5.02%  '
Set-TextValue 'E4' '  +0.10%  '
Set-TextValue 'D5' '549.16'
Set-TextValue 'E5' '  -0.77%  '
Set-TextValue 'D6' '172.75'
Set-TextValue 'E6' '  -3.76%  '
Set-TextValue 'E7' '  -4.93%  '
Set-TextValue 'D9' '3.310.84'
Set-TextValue 'E9' '  -5.03%  '
Set-TextValue 'E10' '  -3.32%  '
Set-TextValue 'E11' '  -2.77%  '
Set-TextValue 'D12' '53.07'
Set-TextValue 'E12' '  -1.69%  '
Set-TextValue 'D13' '0.0000264'
Set-TextValue 'E13' '  -2.93%  '
Set-TextValue 'E14' '  -3.18%  '
Set-TextValue 'D15' '3.845.17'
Set-TextValue 'E15' '  -4.98%  '
Set-TextValue 'D16' '18.17'
Set-TextValue 'E16' '  -1.24%  '
Set-TextValue 'E17' '  -3.41%  '
Set-TextValue 'D18' '3.313.33'
Set-TextValue 'E18' '  -5.00%  '
Set-TextValue 'D19' '11.69'
Set-TextValue 'E19' '  -4.13%  '
Set-TextValue 'D20' '63.418.59'
Set-TextValue 'E20' '  -3.16%  '
Set-TextValue 'D21' '0.967'
Set-TextValue 'E21' '  -2.74%  '
Set-TextValue 'D22' '422.90'
Set-TextValue 'E22' '  +2.04%  '
Set-TextValue 'D23' '4.42'
Set-TextValue 'E23' '  +7.78%  '
Set-TextValue 'E24' '  -0.17%  '
Set-TextValue 'D25' '13.37'
Set-TextValue 'E25' '  +4.63%  '
Set-TextValue 'D26' '83.01'
Set-TextValue 'E26' '  -3.18%  '
Set-TextValue 'D27' '10.61'
Set-TextValue 'E27' '  -1.63%  '
Set-TextValue 'E28' '  -4.43%  '
Set-TextValue 'D29' '8.65'
Set-TextValue 'E29' '  -4.26%  '
Set-TextValue 'D30' '29.15'
Set-TextValue 'E30' '  -3.96%  '
Set-TextValue 'D31' '6.45'
Set-TextValue 'E31' '  -0.34%  '
Set-TextValue 'D32' '11.37'
Set-TextValue 'E32' '  -2.31%  '
Set-TextValue 'D33' '577.40'
Set-TextValue 'E33' '  -6.89%  '
Set-TextValue 'E34' '  -3.59%  '
Set-TextValue 'E35' '  -2.04%  '
Set-TextValue 'E36' '  -0.14%  '
Set-TextValue 'E37' '  -1.18%  '
Set-TextValue 'D38' '3.44'
Set-TextValue 'E38' '  +4.99%  '
Set-TextValue 'D39' '35.08'
Set-TextValue 'E39' '  -5.41%  '
Set-TextValue 'D40' '0.0₃0736'
Set-TextValue 'E40' '  -6.84%  '
Set-TextValue 'D41' '0.364'
Set-TextValue 'E41' '  -4.30%  '
Set-TextValue 'D42' '3.117.48'
Set-TextValue 'E42' '  -7.41%  '
Set-TextValue 'E43' '  +0.01%  '
Set-TextValue 'E44' '  -2.48%  '
Set-TextValue 'D45' '3.16'
Set-TextValue 'E45' '  -2.89%  '
Set-TextValue 'E46' '  -3.19%  '
Set-TextValue 'E47' '  -3.81%  '
Set-TextValue 'E48' '  -6.24%  '
Set-TextValue 'E49' '  -4.09%  '
Set-TextValue 'D50' '135.03'
Set-TextValue 'E50' '  -1.90%  '
Set-TextValue 'E51' '  -4.72%  '
